# Sample Project1 - Main.xlsx save
# Semantic change: Rules!C8 ("Integer min" for rule R10) updated from 0 to 1110.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C8").Value = 1110
